$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-04-22 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-23 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("981÷7=140, 1", $true, $false, $false, $false, $false, $true, 1, $false, "746÷8=93, 2", 2) | Out-Null
$d.Content.Find.Execute("135÷7=19, 2", $true, $false, $false, $false, $false, $true, 1, $false, "985÷3=328, 1", 2) | Out-Null
$d.Content.Find.Execute("387÷8=48, 3", $true, $false, $false, $false, $false, $true, 1, $false, "893÷8=111, 5", 2) | Out-Null
$d.Content.Find.Execute("155÷3=51, 2", $true, $false, $false, $false, $false, $true, 1, $false, "445÷9=49, 4", 2) | Out-Null
$d.Content.Find.Execute("600÷2=300, 0", $true, $false, $false, $false, $false, $true, 1, $false, "612÷8=76, 4", 2) | Out-Null
$d.Content.Find.Execute("261÷2=130, 1", $true, $false, $false, $false, $false, $true, 1, $false, "272÷9=30, 2", 2) | Out-Null
$d.Content.Find.Execute("535÷3=178, 1", $true, $false, $false, $false, $false, $true, 1, $false, "588÷2=294, 0", 2) | Out-Null
$d.Content.Find.Execute("869÷7=124, 1", $true, $false, $false, $false, $false, $true, 1, $false, "385÷5=77, 0", 2) | Out-Null
$d.Content.Find.Execute("249÷6=41, 3", $true, $false, $false, $false, $false, $true, 1, $false, "276÷6=46, 0", 2) | Out-Null
$d.Content.Find.Execute("266÷8=33, 2", $true, $false, $false, $false, $false, $true, 1, $false, "303÷2=151, 1", 2) | Out-Null
$d.Content.Find.Execute("775÷3=258, 1", $true, $false, $false, $false, $false, $true, 1, $false, "550÷3=183, 1", 2) | Out-Null
$d.Content.Find.Execute("555÷7=79, 2", $true, $false, $false, $false, $false, $true, 1, $false, "353÷3=117, 2", 2) | Out-Null
$d.Content.Find.Execute("789÷5=157, 4", $true, $false, $false, $false, $false, $true, 1, $false, "137÷2=68, 1", 2) | Out-Null
$d.Content.Find.Execute("931÷8=116, 3", $true, $false, $false, $false, $false, $true, 1, $false, "851÷7=121, 4", 2) | Out-Null
$d.Content.Find.Execute("503÷4=125, 3", $true, $false, $false, $false, $false, $true, 1, $false, "649÷8=81, 1", 2) | Out-Null
$d.Content.Find.Execute("419÷6=69, 5", $true, $false, $false, $false, $false, $true, 1, $false, "183÷9=20, 3", 2) | Out-Null
$d.Content.Find.Execute("136÷6=22, 4", $true, $false, $false, $false, $false, $true, 1, $false, "916÷6=152, 4", 2) | Out-Null
$d.Content.Find.Execute("995÷6=165, 5", $true, $false, $false, $false, $false, $true, 1, $false, "386÷8=48, 2", 2) | Out-Null
$d.Content.Find.Execute("773÷7=110, 3", $true, $false, $false, $false, $false, $true, 1, $false, "110÷9=12, 2", 2) | Out-Null
$d.Content.Find.Execute("130÷4=32, 2", $true, $false, $false, $false, $false, $true, 1, $false, "126÷5=25, 1", 2) | Out-Null
$d.Content.Find.Execute("429÷2=214, 1", $true, $false, $false, $false, $false, $true, 1, $false, "688÷2=344, 0", 2) | Out-Null
$d.Content.Find.Execute("889÷5=177, 4", $true, $false, $false, $false, $false, $true, 1, $false, "750÷7=107, 1", 2) | Out-Null
$d.Content.Find.Execute("253÷8=31, 5", $true, $false, $false, $false, $false, $true, 1, $false, "876÷6=146, 0", 2) | Out-Null
$d.Content.Find.Execute("740÷7=105, 5", $true, $false, $false, $false, $false, $true, 1, $false, "508÷9=56, 4", 2) | Out-Null
$d.Content.Find.Execute("642÷3=214, 0", $true, $false, $false, $false, $false, $true, 1, $false, "216÷6=36, 0", 2) | Out-Null
